$d = $word.ActiveDocument

$replacements = @(
    @{old = "108×8="; new = "324×4="},
    @{old = "809×5="; new = "748×5="},
    @{old = "378×7="; new = "679×6="},
    @{old = "530×7="; new = "572×4="},
    @{old = "299×7="; new = "707×5="},
    @{old = "421×4="; new = "239×6="},
    @{old = "225×6="; new = "370×3="},
    @{old = "322×5="; new = "529×3="},
    @{old = "889×6="; new = "376×6="},
    @{old = "837×9="; new = "443×4="},
    @{old = "635×6="; new = "194×9="},
    @{old = "336×2="; new = "610×5="},
    @{old = "221×7="; new = "633×3="},
    @{old = "958×9="; new = "670×4="},
    @{old = "227×8="; new = "772×4="},
    @{old = "689×8="; new = "394×2="},
    @{old = "673×6="; new = "431×7="},
    @{old = "522×4="; new = "877×2="},
    @{old = "116×8="; new = "594×2="},
    @{old = "345×2="; new = "676×6="},
    @{old = "537×7="; new = "952×9="},
    @{old = "113×4="; new = "892×8="},
    @{old = "112×2="; new = "697×7="},
    @{old = "679×3="; new = "495×8="},
    @{old = "849×6="; new = "822×4="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
